$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F4").Value = 3
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = 8
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 1
